$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(8616,8616,8504,8504,8369,8369,8369,8369,8369,8369,8094,8094,8094,8094,8094,8094,8094,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,7581,7581,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7534,7320,7320,7320,7320,7320,7320,7320,7320,7320,7320,7318,7318,7318,7310,7310,7310,7310,7310,7310,7310,7310)

$arr = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $arr[$i,0] = $values[$i]
}

$ws.Range("C2:C199").Value = $arr

Write-Output "Updated C2:C199 with new fitness values"
